$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("R2").Value = 1.547445255474453
$ws.Range("AI2").Value = 1.547445255474453
$ws.Range("AZ2").Value = 1.547445255474453
$ws.Range("J4").Value = 2
$ws.Range("AA4").Value = 2
$ws.Range("AR4").Value = 2
$ws.Range("R5").Value = 2.759124087591241
$ws.Range("AI5").Value = 2.759124087591241
$ws.Range("AZ5").Value = 2.759124087591241
$ws.Range("J7").Value = 1.09004583558668
$ws.Range("AA7").Value = 1.09004583558668
$ws.Range("AR7").Value = 1.09004583558668
$ws.Range("J8").Value = 1.09004583558668
$ws.Range("AA8").Value = 1.09004583558668
$ws.Range("AR8").Value = 1.09004583558668
$ws.Range("R11").Value = 1.875912408759124
$ws.Range("AI11").Value = 1.875912408759124
$ws.Range("AZ11").Value = 1.875912408759124
$ws.Range("R19").Value = 2.08029197080292
$ws.Range("AI19").Value = 2.08029197080292
$ws.Range("AZ19").Value = 2.08029197080292
$ws.Range("R20").Value = 1.284671532846715
$ws.Range("AI20").Value = 1.284671532846715
$ws.Range("AZ20").Value = 1.284671532846715
$ws.Range("R33").Value = 1.751824817518248
$ws.Range("AI33").Value = 1.751824817518248
$ws.Range("AZ33").Value = 1.751824817518248
$ws.Range("R34").Value = 2.248175182481752
$ws.Range("AI34").Value = 2.248175182481752
$ws.Range("AZ34").Value = 2.248175182481752
$ws.Range("R37").Value = 2.058394160583942
$ws.Range("AI37").Value = 2.058394160583942
$ws.Range("AZ37").Value = 2.058394160583942
$ws.Range("R38").Value = 2.124087591240876
$ws.Range("AI38").Value = 2.124087591240876
$ws.Range("AZ38").Value = 2.124087591240876
$ws.Range("I44").Value = 2.983213820359353
$ws.Range("Z44").Value = 2.983213820359353
$ws.Range("AQ44").Value = 2.983213820359353
$ws.Range("R45").Value = 2.233576642335767
$ws.Range("AI45").Value = 2.233576642335767
$ws.Range("AZ45").Value = 2.233576642335767
$ws.Range("R46").Value = 2.35036496350365
$ws.Range("AI46").Value = 2.35036496350365
$ws.Range("AZ46").Value = 2.35036496350365
$ws.Range("R50").Value = 1.547445255474453
$ws.Range("AI50").Value = 1.547445255474453
$ws.Range("AZ50").Value = 1.547445255474453
$ws.Range("R52").Value = 1.992700729927007
$ws.Range("AI52").Value = 1.992700729927007
$ws.Range("AZ52").Value = 1.992700729927007
$ws.Range("R61").Value = 2.029197080291971
$ws.Range("AI61").Value = 2.029197080291971
$ws.Range("AZ61").Value = 2.029197080291971
$ws.Range("R62").Value = 2.248175182481752
$ws.Range("AI62").Value = 2.248175182481752
$ws.Range("AZ62").Value = 2.248175182481752
$ws.Range("R64").Value = 1.445255474452555
$ws.Range("AI64").Value = 1.445255474452555
$ws.Range("AZ64").Value = 1.445255474452555
$ws.Range("R74").Value = 2.248175182481752
$ws.Range("AI74").Value = 2.248175182481752
$ws.Range("AZ74").Value = 2.248175182481752
$ws.Range("I77").Value = 1
$ws.Range("Z77").Value = 1
$ws.Range("AQ77").Value = 1
$ws.Range("J80").Value = 1.13506875338002
$ws.Range("AA80").Value = 1.13506875338002
$ws.Range("AR80").Value = 1.13506875338002
$ws.Range("R83").Value = 4.072992700729927
$ws.Range("AI83").Value = 4.072992700729927
$ws.Range("AZ83").Value = 4.072992700729927
$ws.Range("J98").Value = 1.136529288015673
$ws.Range("AA98").Value = 1.136529288015673
$ws.Range("AR98").Value = 1.136529288015673
$ws.Range("I101").Value = 1.355945793324701
$ws.Range("Z101").Value = 1.355945793324701
$ws.Range("AQ101").Value = 1.355945793324701
$ws.Range("J107").Value = 2
$ws.Range("AA107").Value = 2
$ws.Range("AR107").Value = 2
$ws.Range("R108").Value = 2.124087591240876
$ws.Range("AI108").Value = 2.124087591240876
$ws.Range("AZ108").Value = 2.124087591240876
$ws.Range("R109").Value = 3.489051094890511
$ws.Range("AI109").Value = 3.489051094890511
$ws.Range("AZ109").Value = 3.489051094890511
$ws.Range("I110").Value = 1.957688220192625
$ws.Range("Z110").Value = 1.957688220192625
$ws.Range("AQ110").Value = 1.957688220192625
$ws.Range("J115").Value = 1.09004583558668
$ws.Range("AA115").Value = 1.09004583558668
$ws.Range("AR115").Value = 1.09004583558668
$ws.Range("I118").Value = 1.82639681360069
$ws.Range("R118").Value = 1.700729927007299
$ws.Range("Z118").Value = 1.82639681360069
$ws.Range("AI118").Value = 1.700729927007299
$ws.Range("AQ118").Value = 1.82639681360069
$ws.Range("AZ118").Value = 1.700729927007299
$ws.Range("R120").Value = 1.700729927007299
$ws.Range("AI120").Value = 1.700729927007299
$ws.Range("AZ120").Value = 1.700729927007299
$ws.Range("R121").Value = 1.708029197080292
$ws.Range("AI121").Value = 1.708029197080292
$ws.Range("AZ121").Value = 1.708029197080292
$ws.Range("R124").Value = 2.08029197080292
$ws.Range("AI124").Value = 2.08029197080292
$ws.Range("AZ124").Value = 2.08029197080292
$ws.Range("R132").Value = 1.708029197080292
$ws.Range("AI132").Value = 1.708029197080292
$ws.Range("AZ132").Value = 1.708029197080292
$ws.Range("R134").Value = 2.124087591240876
$ws.Range("AI134").Value = 2.124087591240876
$ws.Range("AZ134").Value = 2.124087591240876
$ws.Range("R139").Value = 1.708029197080292
$ws.Range("AI139").Value = 1.708029197080292
$ws.Range("AZ139").Value = 1.708029197080292
$ws.Range("R160").Value = 2.248175182481752
$ws.Range("AI160").Value = 2.248175182481752
$ws.Range("AZ160").Value = 2.248175182481752
$ws.Range("R168").Value = 1.700729927007299
$ws.Range("AI168").Value = 1.700729927007299
$ws.Range("AZ168").Value = 1.700729927007299
$ws.Range("J178").Value = 1.09004583558668
$ws.Range("AA178").Value = 1.09004583558668
$ws.Range("AR178").Value = 1.09004583558668
$ws.Range("J179").Value = 1.13506875338002
$ws.Range("AA179").Value = 1.13506875338002
$ws.Range("AR179").Value = 1.13506875338002
$ws.Range("R197").Value = 2.233576642335767
$ws.Range("AI197").Value = 2.233576642335767
$ws.Range("AZ197").Value = 2.233576642335767
$ws.Range("R199").Value = 2.277372262773723
$ws.Range("AI199").Value = 2.277372262773723
$ws.Range("AZ199").Value = 2.277372262773723
